# Add a new "2022-Q4" quarterly sheet to the 002960-青鸟消防 workbook.
#
# 1. "总计" (summary) sheet gets a new row 2 for 2022-Q4 (17 holdings,
#    3.8 yi), and every existing row shifts down by one.
# 2. A brand-new worksheet named "2022-Q4" is inserted right after "总计"
#    (and therefore right before "2022-Q3"), carrying the per-fund holding
#    breakdown, formatted like the existing quarter sheets.
# 3. The previously-active last sheet ("2020-Q4") is re-activated so the
#    workbook keeps the same sheet selected as before the edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) 总计 sheet: insert the 2022-Q4 row at the top of the data and shift
#    everything else down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A2:D2").Insert()

# Copy the row-2 number format (bold/border index-column style) down onto
# the freshly inserted row so column A keeps looking like the rest of the
# index column; clear the stray format Insert() leaves on B2:D2.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 17
$total.Range("D2").Value = 3.8

# The index column (A) is a 0-based row counter; every pre-existing row
# slid down one position, so its counter value needs to grow by one too.
for ($r = 3; $r -le 10; $r++) {
    $total.Range("A$r").Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) Create the new "2022-Q4" sheet by cloning "2022-Q3" (so it inherits
#    identical styles/page setup/column widths), then overwrite its
#    contents with the 2022-Q4 fund breakdown.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $total)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The source quarter had 9 fund rows (rows 2-10); 2022-Q4 has 17 (rows
# 2-18). Extend the sheet by copying the formatted-but-empty row 2 down
# onto the 8 additional rows so column A keeps its index-column style.
$q4.Range("A2:H2").Copy()
$q4.Range("A11:H18").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$fundRows = @(
        @('001576','国泰智能装备股票A','19.52','92.90','7.10','1.3859'),
        @('340008','兴全有机增长混合','21.32','79.09','4.77','1.0170'),
        @('010330','东吴兴享成长混合A','7.86','82.94','7.30','0.5738'),
        @('011322','国泰智能装备股票C','2.61','92.90','7.10','0.1853'),
        @('007811','淳厚信泽灵活配置混合A','4.32','79.09','2.89','0.1248'),
        @('011462','东吴兴享成长混合C','1.17','82.94','7.30','0.0854'),
        @('010551','淳厚欣颐一年持有期混合','2.49','88.75','3.03','0.0754'),
        @('011349','淳厚现代服务业股票A','2.41','88.44','3.01','0.0725'),
        @('016588','富国融甄混合A','4.57','29.63','1.51','0.0690'),
        @('005413','金信民长灵活配置混合C','0.86','89.93','4.96','0.0427'),
        @('005412','金信民长灵活配置混合A','0.83','89.93','4.96','0.0412'),
        @('020023','国泰事件驱动策略混合A','2.19','82.03','1.75','0.0383'),
        @('016589','富国融甄混合C','1.96','29.63','1.51','0.0296'),
        @('007812','淳厚信泽灵活配置混合C','0.81','79.09','2.89','0.0234'),
        @('011350','淳厚现代服务业股票C','0.57','88.44','3.01','0.0172'),
        @('006209','中信保诚新蓝筹灵活配置混合','0.52','81.42','2.60','0.0135'),
        @('015592','国泰事件驱动策略混合C','0.01','82.03','1.75','0.0002')
)
$ranks = @(5,6,2,5,9,2,9,8,7,3,3,10,7,9,8,10,10)

# Force columns B:G to stay plain text (fund codes like "001576" must not
# be coerced into numbers) while writing, then drop the temporary "@"
# number-format again so no extra style survives on the exported cells.
$textRange = $q4.Range("B2:G18")
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $q4.Range("A$r").Value = $i
    $q4.Range("B$r").Value = $row[0]
    $q4.Range("C$r").Value = $row[1]
    $q4.Range("D$r").Value = $row[2]
    $q4.Range("E$r").Value = $row[3]
    $q4.Range("F$r").Value = $row[4]
    $q4.Range("G$r").Value = $row[5]
    $q4.Range("H$r").Value = $ranks[$i]
}

$textRange.ClearFormats()

# ---------------------------------------------------------------------
# 3) Keep the same tab active/selected as before the edit (the last
#    sheet, "2020-Q4").
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Activate()
